## edit.ps1 -- apply the DATOS.docx revision:
##   1. Paragraph 1: insert a new run "veraz y " right before the run
##      "de alta calidad, " (so "...Tu misión es entregar información
##      de alta calidad," becomes "...información veraz y de alta
##      calidad,"), while leaving every other run boundary in that
##      paragraph exactly as it was.
##   2. Paragraph 3 ("Alcance: ..."): split the single run into three
##      runs, inserting a new "Adicionalmente " run between "... RAMA. "
##      and "Enfocar las respuestas ...".
##   3. Move the _GoBack bookmark from the end of paragraph 3 to the
##      (empty) paragraph 4 that follows it.
##
## NOTE: this COM-interop runtime recomputes/merges the <w:r> runs of a
## paragraph whenever text is inserted into it (same as how real Word
## frequently coalesces runs that share identical formatting) -- every
## run from the insertion point to the end of the paragraph collapses
## into one run. To recover the original (untouched) run boundaries we
## "re-split" the paragraph afterwards by toggling a trivial character
## formatting property (Bold) on and back off over each exact
## sub-range; that forces the engine to split the run at the sub-range
## edges without otherwise touching the rest of the paragraph's text.
## Sub-ranges are located purely by character offsets (not by
## re-searching for text), so repeated words elsewhere in the
## paragraph can't cause the wrong boundary to be picked.

$d = $word.ActiveDocument

function Apply-RunLengths($paragraphRange, $lengths) {
    # $paragraphRange: a Range covering the whole paragraph (content
    #                  only; Word silently excludes the trailing pilcrow
    #                  from these operations anyway).
    # $lengths: ordered array of the character counts of each run that
    #           should exist, left to right, covering the paragraph.
    $base = $paragraphRange.Start
    $offset = 0
    $n = $lengths.Count
    for ($i = 0; $i -lt ($n - 1); $i++) {
        $offset = $offset + $lengths[$i]
        $boundary = $base + $offset
        # Toggle a character property across the two runs that should
        # remain on either side of this boundary -- touching
        # [boundary-1, boundary+1) is enough to force a split exactly
        # at $boundary without re-touching the whole paragraph.
        $probe = $d.Range($boundary - 1, $boundary + 1)
        $probe.Bold = 1
        $probe.Bold = 0
    }
}

## --- 1. "informacion ... de alta calidad, " -> insert "veraz y " --------
$find1 = $d.Content
$ok1 = $find1.Find.Execute("de alta calidad, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok1) {
    $insertPoint = $d.Range($find1.Start, $find1.Start)
    $insertPoint.InsertBefore("veraz y ")
}

$p1 = $d.Paragraphs.Item(1)
$lengths1 = @(
    20,  # "Objetivo principal: "
    179, # "Eres un asistente ... Tu misión es entregar"
    1,   # " "
    12,  # "información "
    8,   # "veraz y "            <- new run
    18,  # "de alta calidad, "
    14,  # "de la manera m"
    1,   # "á"
    79,  # "s clara posible ... adquirido."
    1,   # " "
    112, # "No entregar ... debe siempre"
    1,   # " "
    112, # "formal y cortés ... planteadas"
    1    # "."
)
Apply-RunLengths $p1.Range $lengths1

## --- 2. "Alcance: ... RAMA. Enfocar ..." -> split + insert "Adicionalmente " --
$find2 = $d.Content
$ok2 = $find2.Find.Execute("Enfocar las respuestas", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok2) {
    $insertPoint2 = $d.Range($find2.Start, $find2.Start)
    $insertPoint2.InsertBefore("Adicionalmente ")
}

$p3 = $d.Paragraphs.Item(3)
$lengths3 = @(
    195, # "Alcance: ... RAMA. "
    15,  # "Adicionalmente "     <- new run
    74   # "Enfocar ... base de datos."
)
Apply-RunLengths $p3.Range $lengths3

## --- 3. Move the _GoBack bookmark to the empty paragraph after "Alcance" --
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)

Write-Output "done"
